$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 6, shifting existing rows (and the
# "Melbourne Airport" tail rows at 121-123) down by one, extending the
# table through row 124.
$ws.Rows.Item(6).Insert()

$ws.Range("A6").Value = "12-18 Distribution Dr, Truganina VIC 3029"
$ws.Range("B6").Value = -37.814598
$ws.Range("C6").Value = 144.762011
$ws.Range("D6").Value = "Wyndham (C)"
